$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Rows(29).Insert()
$ws.Range("A29").Value = "OCNet(Res101) (Yuan et al. 2018)"
$ws.Range("B29").Value = 81.7
$ws.Range("C29").HorizontalAlignment = -4108
